# Natmi following Dr Hou advice
# Update Ligand/Receptor-expressing cell counts (E, K) from 1 to 3 for all
# data rows, and refresh the dependent expression/specificity statistics
# (G,H,I,J,M,N,O,P,Q,R,S,T) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by letter) that change, in order, for each data row 2..10:
# E, G, H, I, J, K, M, N, O, P, Q, R, S, T
$cols = @("E","G","H","I","J","K","M","N","O","P","Q","R","S","T")

# Row number followed by the 14 new values in the order of $cols above.
$data = @(
    @(2,  3, 21.870458, 65.61137400000001, 0.9662106525726075, 0.9662106525726075, 3, 91.47356466666666, 274.420694, 0.7914295280104694, 0.7914295280104694, 2000.568754152617, 18005.11878737356, 0.7646876407242263, 0.7646876407242263),
    @(3,  3, 21.870458, 65.61137400000001, 0.9662106525726075, 0.9662106525726075, 3, 20.94207233333333, 62.82621699999999, 0.181190866265331, 0.181190866265331, 458.0127133991287, 4122.114420592158, 0.1750685451344215, 0.1750685451344215),
    @(4,  3, 21.870458, 65.61137400000001, 0.9662106525726075, 0.9662106525726075, 3, 3.164539666666667, 9.493619000000001, 0.02737960572419959, 0.02737960572419959, 69.20993186916735, 622.8893868225061, 0.02645446671395959, 0.02645446671395959),
    @(5,  3, 0.1926236666666667, 0.577871, 0.008509882997005752, 0.008509882997005752, 3, 91.47356466666666, 274.420694, 0.7914295280104694, 0.7914295280104694, 17.61997342916378, 158.579760862474, 0.006734972683744581, 0.006734972683744581),
    @(6,  3, 0.1926236666666667, 0.577871, 0.008509882997005752, 0.008509882997005752, 3, 20.94207233333333, 62.82621699999999, 0.181190866265331, 0.181190866265331, 4.033938760445222, 36.305448844007, 0.001541913072044084, 0.001541913072044084),
    @(7,  3, 0.1926236666666667, 0.577871, 0.008509882997005752, 0.008509882997005752, 3, 3.164539666666667, 9.493619000000001, 0.02737960572419959, 0.02737960572419959, 0.6095652339054445, 5.486087105149001, 0.0002329972412170875, 0.0002329972412170875),
    @(8,  3, 0.572208, 1.716624, 0.02527946443038672, 0.02527946443038671, 3, 91.47356466666666, 274.420694, 0.7914295280104694, 0.7914295280104694, 52.341905490784, 471.077149417056, 0.02000691460249841, 0.02000691460249841),
    @(9,  3, 0.572208, 1.716624, 0.02527946443038672, 0.02527946443038671, 3, 20.94207233333333, 62.82621699999999, 0.181190866265331, 0.181190866265331, 11.983221325712, 107.848991931408, 0.004580408058865392, 0.004580408058865391),
    @(10, 3, 0.572208, 1.716624, 0.02527946443038672, 0.02527946443038671, 3, 3.164539666666667, 9.493619000000001, 0.02737960572419959, 0.02737960572419959, 1.810774913584, 16.296974222256, 0.0006921417690229162, 0.0006921417690229161)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}
